# "Tried to implement Penality Reward System (unfinished)"
#
# Weekly Quantity sheet: collapse the long weekly tail into a shortened
# "reward/penalty" style window - the last two weekly rows move up and get
# new dates/quantities, and the remaining trailing rows are removed.
#
# Monthly Trend sheet: same idea - the last monthly row moves up and gets
# a new date/quantity, and the remaining trailing rows are removed.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" (A1:B16 -> A1:B9) ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

$ws1.Range("A8").Value = 45144.99999999999
$ws1.Range("B8").Value = 6

$ws1.Range("A9").Value = 45151.99999999999
$ws1.Range("B9").Value = 6

# Remove the now-redundant trailing rows 10-16.
$ws1.Range("A10:B16").EntireRow.Delete()

# --- Sheet 2: "Monthly Trend" (A1:B7 -> A1:B5) ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Range("A5").Value = 45169.99999999999
$ws2.Range("B5").Value = 12

# Remove the now-redundant trailing rows 6-7.
$ws2.Range("A6:B7").EntireRow.Delete()
